$d = $word.ActiveDocument

# Locate the paragraph that holds the field (fldChar/instrText "{ m:false.yesNo() }"):
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the field to rewrite."
}

# Replace the field-code runs (fldChar begin/instrText.../fldChar end) with
# plain literal-text runs spelling out the same field code, keeping the
# bookmark and the coloured "false" run intact. This mirrors the parser
# switch to TokenIteratorFieldRewriterSplit, which emits split literal
# runs instead of a real Word field.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979"><w:r><w:t>{</w:t></w:r><w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r w:rsidR="00FA1F93"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>false</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="00FA1F93"><w:t>yesNo</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="005E0B49"><w:t>()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($xml)
